# Apply the Tutorial 6 solution update:
#  - Replace slashes with dashes in the date strings of column A (rows 3-21)
#  - Update the attendance-count columns (D, E, G, H) for a handful of rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: change date separator from "/" to "-" for rows 3..21 ---
# Force the cells to remain plain text (some "dd-mm-yyyy" strings with
# dd <= 12 would otherwise be auto-parsed as a date serial by Excel),
# matching the original inline-string representation.
$ws.Range("A3:A21").NumberFormat = "@"

$ws.Range("A3").Value = "28-07-2022"
$ws.Range("A4").Value = "01-08-2022"
$ws.Range("A5").Value = "04-08-2022"
$ws.Range("A6").Value = "08-08-2022"
$ws.Range("A7").Value = "11-08-2022"
$ws.Range("A8").Value = "15-08-2022"
$ws.Range("A9").Value = "18-08-2022"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("A13").Value = "01-09-2022"
$ws.Range("A14").Value = "05-09-2022"
$ws.Range("A15").Value = "08-09-2022"
$ws.Range("A16").Value = "12-09-2022"
$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A21").Value = "29-09-2022"

# --- Attendance value updates (D, E, G, H columns) ---
# Row 3: Total=1, Invalid=1
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4: Total=1, Real=1, Absent=0
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Row 5: Total=1, Real=1, Absent=0
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# Row 6: Total=1, Real=1, Absent=0
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0

# Row 12: Total=1, Invalid=1
$ws.Range("D12").Value = 1
$ws.Range("G12").Value = 1
